$wb = $excel.ActiveWorkbook

# Sheet "建物" (building): property_category column (I) was "land" for every
# row; these rows are actually buildings, so set it to "building".
$wsBuilding = $wb.Worksheets.Item("建物")
for ($r = 2; $r -le 7; $r++) {
    $wsBuilding.Range("I" + $r).Value = "building"
}

# Sheet "汽車" (car): property_category column (H) was "land"; this row is
# actually a car, so set it to "car".
$wsCar = $wb.Worksheets.Item("汽車")
$wsCar.Range("H2").Value = "car"
